# Update rules in DiscountRules.xlsx
# - Adds a new (5th) column "E" to the existing rule table, with header
#   "ACTION" on the NAME/CONDITION/ACTION row and "Test" on the snippet
#   row beneath it; all the other existing rows simply gain a blank
#   cell in column E so the sheet's used range grows to column E.
# - Appends a new rule row (row 24) to the decision table with an empty
#   rule name and a new action value "5x5 0048" in the new column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E: create (empty) cells for all the pre-existing rows ---
# Rows 1-17 are the decision table's property/metadata rows.
# Rows 20-23 are the existing rule rows below the table header/snippet.
# Using a formatting no-op (border style) forces Excel to materialize
# the cell without altering its (default) value/style.
$ws.Range("E1:E17").Borders.LineStyle = -4142
$ws.Range("E20:E23").Borders.LineStyle = -4142

# --- Column E header + snippet row for the rule table ---
$ws.Cells.Item(18, 5).Value = "ACTION"
$ws.Cells.Item(19, 5).Value = "Test"

# --- New rule row 24 ---
# Columns A-D stay blank (new/empty rule name, no conditions/action set)
$ws.Range("A24:D24").Borders.LineStyle = -4142
# New action value for the new "5x5 0048" rule/column
$ws.Cells.Item(24, 5).Value = "5x5 0048"
